# Amended a couple of entries in the risk assessment
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- RSI risk (row 12): broaden the response to also mention wrist supports ---
$ws.Range("D12").Value = "Ensure I move wrists about and wear wrist supports"

# --- Replace the "Global Pandemic" risk (row 18) with a new "Possible World War 3" risk ---
$ws.Range("B18").Value = "Possible World War 3"
$ws.Range("C18").Value = "A large nation could invade its neighbour with the potential to spark a nuclear conflict between world superpowers"
$ws.Range("D18").Value = 'Read "When the Wind Blows" to relax'
# Response / Likelihood / Impact (E18:G18) stay as they were:
#   "Just keep going, one way or another" / "Super Unlikely" / "Still working that out"

# --- Update the view so the window is scrolled down and G18 is the active selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("G18").Select()
